# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across all
# eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with freshly
# pulled values. Plain data cells, no formulas involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4132.75
$ws.Range("I32").Value = 3987.0908
$ws.Range("K32").Value = 3987.0908
$ws.Range("M32").Value = -3661.0908
$ws.Range("H40").Value = 3499.75
$ws.Range("I40").Value = 2599.6
$ws.Range("K40").Value = 2599.6
$ws.Range("M40").Value = -2424.6
$ws.Range("H94").Value = 50126224
$ws.Range("I94").Value = 100002960
$ws.Range("J94").Value = 249490.8
$ws.Range("K94").Value = 100002960
$ws.Range("L94").Value = 249490.8
$ws.Range("M94").Value = -100002509
$ws.Range("N94").Value = -250392.8
$ws.Range("H113").Value = 17566
$ws.Range("J113").Value = 6349
$ws.Range("L113").Value = 6349
$ws.Range("N113").Value = -12857
$ws.Range("H135").Value = 6773.875
$ws.Range("I135").Value = 8031.8335
$ws.Range("K135").Value = 72286.5015
$ws.Range("M135").Value = -69751.5015
$ws.Range("H137").Value = 258154.66
$ws.Range("I137").Value = 439543.4
$ws.Range("J137").Value = 2076.4119
$ws.Range("K137").Value = 1318630.2
$ws.Range("L137").Value = 6229.2357
$ws.Range("M137").Value = -1316080.2
$ws.Range("N137").Value = -11329.2357
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 36999
$ws.Range("J43").Value = 45499
$ws.Range("L43").Value = 45499
$ws.Range("N43").Value = -46125
$ws.Range("H122").Value = 487536.53
$ws.Range("I122").Value = 3766.3845
$ws.Range("J122").Value = 880599.75
$ws.Range("K122").Value = 11299.1535
$ws.Range("L122").Value = 2641799.25
$ws.Range("M122").Value = -8849.1535
$ws.Range("N122").Value = -2646699.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 647.8
$ws.Range("I22").Value = 647.8
$ws.Range("K22").Value = 647.8
$ws.Range("M22").Value = -474.8
$ws.Range("H86").Value = 5798
$ws.Range("I86").Value = 8258.6
$ws.Range("K86").Value = 8258.6
$ws.Range("M86").Value = -7135.6
$ws.Range("H88").Value = 28409.285
$ws.Range("J88").Value = 28409.285
$ws.Range("L88").Value = 28409.285
$ws.Range("N88").Value = -29221.285
$ws.Range("H89").Value = 5798
$ws.Range("I89").Value = 8258.6
$ws.Range("K89").Value = 41293
$ws.Range("M89").Value = -35677
$ws.Range("H91").Value = 28409.285
$ws.Range("J91").Value = 28409.285
$ws.Range("L91").Value = 28409.285
$ws.Range("N91").Value = -31217.285
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3410.9412
$ws.Range("I58").Value = 2206.2222
$ws.Range("K58").Value = 2206.2222
$ws.Range("M58").Value = -2003.2222
$ws.Range("H68").Value = 99991
$ws.Range("J68").Value = 99991
$ws.Range("L68").Value = 99991
$ws.Range("N68").Value = -101489
$ws.Range("H71").Value = 99991
$ws.Range("J71").Value = 99991
$ws.Range("L71").Value = 299973
$ws.Range("N71").Value = -307461
$ws.Range("I99").Value = 1282800
$ws.Range("J99").Value = 5932.3335
$ws.Range("K99").Value = 1282800
$ws.Range("L99").Value = 5932.3335
$ws.Range("M99").Value = -1281302
$ws.Range("N99").Value = -8928.333500000001
$ws.Range("H107").Value = 8228.817999999999
$ws.Range("I107").Value = 9795.115
$ws.Range("J107").Value = 2411.1428
$ws.Range("K107").Value = 9795.115
$ws.Range("L107").Value = 2411.1428
$ws.Range("M107").Value = -7875.115
$ws.Range("N107").Value = -6251.1428
$ws.Range("I126").Value = 1282800
$ws.Range("J126").Value = 5932.3335
$ws.Range("K126").Value = 3848400
$ws.Range("L126").Value = 17797.0005
$ws.Range("M126").Value = -3845930
$ws.Range("N126").Value = -22737.0005
$ws.Range("H132").Value = 28278.684
$ws.Range("I132").Value = 3792.5
$ws.Range("J132").Value = 96840
$ws.Range("K132").Value = 11377.5
$ws.Range("L132").Value = 290520
$ws.Range("M132").Value = -8847.5
$ws.Range("N132").Value = -295580
$ws.Range("H136").Value = 3410.9412
$ws.Range("I136").Value = 2206.2222
$ws.Range("K136").Value = 6618.6666
$ws.Range("M136").Value = -4068.6666
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 750
$ws.Range("I3").Value = 750
$ws.Range("K3").Value = 2250
$ws.Range("M3").Value = -2138
$ws.Range("H4").Value = 55508630
$ws.Range("I4").Value = 43121044
$ws.Range("K4").Value = 129363132
$ws.Range("M4").Value = -129363020
$ws.Range("H6").Value = 1321.5
$ws.Range("I6").Value = 1283.1111
$ws.Range("J6").Value = 1667
$ws.Range("K6").Value = 3849.3333
$ws.Range("L6").Value = 5001
$ws.Range("M6").Value = -3736.3333
$ws.Range("N6").Value = -5227
$ws.Range("H108").Value = 8758.200000000001
$ws.Range("I108").Value = 7447.75
$ws.Range("K108").Value = 22343.25
$ws.Range("M108").Value = -19463.25
$ws.Range("H131").Value = 1268.75
$ws.Range("I131").Value = 1028.75
$ws.Range("J131").Value = 1748.75
$ws.Range("K131").Value = 3086.25
$ws.Range("L131").Value = 5246.25
$ws.Range("M131").Value = 1953.75
$ws.Range("N131").Value = -15326.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 8864864
$ws.Range("I11").Value = 14250000
$ws.Range("J11").Value = 2402700
$ws.Range("K11").Value = 14250000
$ws.Range("L11").Value = 2402700
$ws.Range("M11").Value = -14249861
$ws.Range("N11").Value = -2402978
$ws.Range("H12").Value = 9003333
$ws.Range("I12").Value = 2002500
$ws.Range("J12").Value = 12503750
$ws.Range("K12").Value = 2002500
$ws.Range("L12").Value = 12503750
$ws.Range("M12").Value = -2002360
$ws.Range("N12").Value = -12504030
$ws.Range("H62").Value = 46000
$ws.Range("J62").Value = 46000
$ws.Range("L62").Value = 46000
$ws.Range("N62").Value = -47372
$ws.Range("H63").Value = 55000
$ws.Range("J63").Value = 55000
$ws.Range("L63").Value = 55000
$ws.Range("N63").Value = -56372
$ws.Range("H65").Value = 46000
$ws.Range("J65").Value = 46000
$ws.Range("L65").Value = 138000
$ws.Range("N65").Value = -144864
$ws.Range("H66").Value = 55000
$ws.Range("J66").Value = 55000
$ws.Range("L66").Value = 165000
$ws.Range("N66").Value = -171864
$ws.Range("H70").Value = 7226.857
$ws.Range("I70").Value = 4894.5
$ws.Range("J70").Value = 8159.8
$ws.Range("K70").Value = 4894.5
$ws.Range("L70").Value = 8159.8
$ws.Range("M70").Value = -4624.5
$ws.Range("N70").Value = -8699.799999999999
$ws.Range("H73").Value = 7226.857
$ws.Range("I73").Value = 4894.5
$ws.Range("J73").Value = 8159.8
$ws.Range("K73").Value = 4894.5
$ws.Range("L73").Value = 8159.8
$ws.Range("M73").Value = -3958.5
$ws.Range("N73").Value = -10031.8
$ws.Range("H132").Value = 3876.4443
$ws.Range("I132").Value = 3810.353
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 11431.059
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -8901.059000000001
$ws.Range("N132").Value = -20060
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4661.077
$ws.Range("J46").Value = 6611.875
$ws.Range("L46").Value = 6611.875
$ws.Range("N46").Value = -6987.875
$ws.Range("H61").Value = 7953.3335
$ws.Range("I61").Value = 5851
$ws.Range("J61").Value = 50000
$ws.Range("K61").Value = 5851
$ws.Range("L61").Value = 50000
$ws.Range("M61").Value = -5649
$ws.Range("N61").Value = -50404
$ws.Range("H93").Value = 7782
$ws.Range("I93").Value = 8455.352999999999
$ws.Range("K93").Value = 8455.352999999999
$ws.Range("M93").Value = -7207.352999999999
$ws.Range("H113").Value = 7953.3335
$ws.Range("I113").Value = 5851
$ws.Range("J113").Value = 50000
$ws.Range("K113").Value = 5851
$ws.Range("L113").Value = 50000
$ws.Range("M113").Value = -3681
$ws.Range("N113").Value = -54340
$ws.Range("H136").Value = 4239.9
$ws.Range("I136").Value = 2012.762
$ws.Range("J136").Value = 6701.4736
$ws.Range("K136").Value = 6038.286
$ws.Range("L136").Value = 20104.4208
$ws.Range("M136").Value = -3488.286
$ws.Range("N136").Value = -25204.4208
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 114246.04
$ws.Range("I62").Value = 209066.27
$ws.Range("J62").Value = 4838.077
$ws.Range("K62").Value = 209066.27
$ws.Range("L62").Value = 4838.077
$ws.Range("M62").Value = -208442.27
$ws.Range("N62").Value = -6086.077
$ws.Range("H65").Value = 114246.04
$ws.Range("I65").Value = 209066.27
$ws.Range("J65").Value = 4838.077
$ws.Range("K65").Value = 1045331.35
$ws.Range("L65").Value = 24190.385
$ws.Range("M65").Value = -1042211.35
$ws.Range("N65").Value = -30430.385
$ws.Range("H82").Value = 42942.25
$ws.Range("J82").Value = 42256.332
$ws.Range("L82").Value = 42256.332
$ws.Range("N82").Value = -43022.332
$ws.Range("H85").Value = 42942.25
$ws.Range("J85").Value = 42256.332
$ws.Range("L85").Value = 42256.332
$ws.Range("N85").Value = -44908.332
$ws.Range("H107").Value = 53498.832
$ws.Range("I107").Value = 4198.8
$ws.Range("K107").Value = 12596.4
$ws.Range("M107").Value = -10676.4
$ws.Range("H113").Value = 5659
$ws.Range("I113").Value = 1788.8
$ws.Range("J113").Value = 10496.75
$ws.Range("K113").Value = 5366.4
$ws.Range("L113").Value = 31490.25
$ws.Range("M113").Value = -3196.4
$ws.Range("N113").Value = -35830.25
$ws.Range("H132").Value = 20511
$ws.Range("I132").Value = 23659.934
$ws.Range("K132").Value = 70979.802
$ws.Range("M132").Value = -68449.802
$ws.Range("H136").Value = 2144.5
$ws.Range("I136").Value = 1736.125
$ws.Range("J136").Value = 2797.9
$ws.Range("K136").Value = 5208.375
$ws.Range("L136").Value = 8393.700000000001
$ws.Range("M136").Value = -2658.375
$ws.Range("N136").Value = -13493.7
